# Clear the contents of A2 (was "-"), leaving it blank, as in the target
# workbook output. All other cells (headers in row 1, numeric data in
# rows 2-4) are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = ""
